$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 previously held a standalone, bordered/bold numeric 0 cell; row 2
# held the question text (as a shared string). Delete row 1 so row 2 (the
# question text) shifts up and becomes the new, only row.
$ws.Rows(1).Delete()

# Reformat the question text in place (Python-dict-repr -> pretty JSON):
#   questions = [{'title': '...', 'ques_type': None, 'options': [], 'score': None}]
# becomes:
#   questions = [
#       {
#           "title": "...",
#           "ques_type": null,
#           "options": [],
#           "score": null
#       }
#   ]
# Using Replace() (rather than re-assigning .Value wholesale) keeps the edit
# a pure text substitution of the existing shared string.
$find1 = "questions = [{'title': '"
$repl1 = "questions = [`n    {`n        ""title"": """
[void]$ws.Cells.Replace($find1, $repl1)

$find2 = "', 'ques_type': None, 'options': [], 'score': None}]"
$repl2 = """,`n        ""ques_type"": null,`n        ""options"": [],`n        ""score"": null`n    }`n]"
[void]$ws.Cells.Replace($find2, $repl2)

# The text now contains embedded newlines, which auto-expands the row
# height; re-run AutoFit so the row settles back to the sheet's standard
# height instead of keeping a stale custom height.
$ws.Rows(1).AutoFit()
